{"js": "const replacements = [\n  [\"2024-03-10 Sunday\", \"2024-03-11 Monday\"],\n  [\"36\u00d730=\", \"62\u00d748=\"],\n  [\"21\u00d747=\", \"87\u00d734=\"],\n  [\"14\u00d715=\", \"66\u00d720=\"],\n  [\"81\u00d753=\", \"76\u00d738=\"],\n  [\"49\u00d738=\", \"51\u00d762=\"],\n  [\"99\u00d760=\", \"33\u00d726=\"],\n  [\"88\u00d712=\", \"13\u00d783=\"],\n  [\"28\u00d757=\", \"38\u00d790=\"],\n  [\"39\u00d788=\", \"53\u00d743=\"],\n  [\"37\u00d753=\", \"37\u00d718=\"],\n  [\"71\u00d796=\", \"41\u00d783=\"],\n  [\"69\u00d796=\", \"55\u00d782=\"],\n  [\"63\u00d722=\", \"98\u00d765=\"],\n  [\"71\u00d793=\", \"47\u00d784=\"],\n  [\"89\u00d736=\", \"91\u00d730=\"],\n  [\"39\u00d776=\", \"68\u00d795=\"],\n  [\"70\u00d770=\", \"17\u00d788=\"],\n  [\"18\u00d782=\", \"72\u00d740=\"],\n  [\"36\u00d786=\", \"91\u00d766=\"],\n  [\"83\u00d778=\", \"64\u00d787=\"],\n  [\"49\u00d754=\", \"51\u00d773=\"],\n  [\"12\u00d780=\", \"11\u00d772=\"],\n  [\"76\u00d731=\", \"44\u00d768=\"],\n  [\"66\u00d722=\", \"68\u00d731=\"],\n  [\"80\u00d767=\", \"46\u00d751=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-03-10 Sunday\", \"2024-03-11 Monday\"),\n    @(\"36\u00d730=\", \"62\u00d748=\"),\n    @(\"21\u00d747=\", \"87\u00d734=\"),\n    @(\"14\u00d715=\", \"66\u00d720=\"),\n    @(\"81\u00d753=\", \"76\u00d738=\"),\n    @(\"49\u00d738=\", \"51\u00d762=\"),\n    @(\"99\u00d760=\", \"33\u00d726=\"),\n    @(\"88\u00d712=\", \"13\u00d783=\"),\n    @(\"28\u00d757=\", \"38\u00d790=\"),\n    @(\"39\u00d788=\", \"53\u00d743=\"),\n    @(\"37\u00d753=\", \"37\u00d718=\"),\n    @(\"71\u00d796=\", \"41\u00d783=\"),\n    @(\"69\u00d796=\", \"55\u00d782=\"),\n    @(\"63\u00d722=\", \"98\u00d765=\"),\n    @(\"71\u00d793=\", \"47\u00d784=\"),\n    @(\"89\u00d736=\", \"91\u00d730=\"),\n    @(\"39\u00d776=\", \"68\u00d795=\"),\n    @(\"70\u00d770=\", \"17\u00d788=\"),\n    @(\"18\u00d782=\", \"72\u00d740=\"),\n    @(\"36\u00d786=\", \"91\u00d766=\"),\n    @(\"83\u00d778=\", \"64\u00d787=\"),\n    @(\"49\u00d754=\", \"51\u00d773=\"),\n    @(\"12\u00d780=\", \"11\u00d772=\"),\n    @(\"76\u00d731=\", \"44\u00d768=\"),\n    @(\"66\u00d722=\", \"68\u00d731=\"),\n    @(\"80\u00d767=\", \"46\u00d751=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
